$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column B. This shifts the existing
# B -> D ("Title of report") and C -> E (blank helper column), and the
# formerly-default D/E -> F/G.
$ws.Columns("B:C").Insert()

# Populate the two new header cells. Set C3 ("Database") before B3
# ("Author") so the shared-string table ends up in the same order as the
# target workbook (index 1 = Database, index 2 = Author).
$ws.Range("C3").Value = "Database"
$ws.Range("B3").Value = "Author"

# Copy the formatting (fill/font/wrap) of the banner row from the old
# title cell (now D1) onto the two new banner cells B1/C1.
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$null = $excel.CutCopyMode

# Copy the formatting of the title cell (now D3) onto the two new header
# cells B3/C3, without disturbing the text we just put in them.
$ws.Range("D3").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$null = $excel.CutCopyMode

# Resize the four data columns to their final widths.
$ws.Columns("B").ColumnWidth = 20.666666666666668
$ws.Columns("C").ColumnWidth = 24.333333333333332
$ws.Columns("D").ColumnWidth = 64.66666666666667
$ws.Columns("E").ColumnWidth = 16.333333333333332

# The "Title" defined name pointed at the old location of the title cell
# (B3); repoint it at the title cell's new address (D3).
$titleName = $wb.Names.Item("Title")
$titleName.RefersTo = "=Sheet1!`$D`$3"

# Match the saved selection/active cell.
$null = $ws.Range("D4").Select()
